# "Add Post Meeting Checkpoint"
# Slide 7 has two summary text boxes whose headline group-name run needs
# updating:
#   Textfeld 3 (shape 1): "GESIS"  -> "Allensbach"
#   Textfeld 4 (shape 2): "GBS"    -> "GESIS"
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# --- Shape 1 ("Textfeld 3"): first run "GESIS" -> "Allensbach" ---
$shp1 = $s.Shapes.Item(1)
$tr1 = $shp1.TextFrame.TextRange

# Remove the old "GESIS" (chars 1-5) outright, then insert the new word
# in front of the following run (the ":" run) and re-apply Bold to just
# the inserted span. Doing it this way (delete, then insert+restyle)
# yields a clean <a:rPr> for the new run instead of carrying over the
# old run's stale spell-check "err" flag the way a plain text
# replacement would.
$old1 = $tr1.Characters(1, 5)
$old1.Text = ""
$afterDelete1 = $tr1.Characters(1, 1)
$afterDelete1.InsertBefore("Allensbach") | Out-Null
$newWord1 = $tr1.Characters(1, 10)
$newWord1.Font.Bold = $true

# --- Shape 2 ("Textfeld 4"): first run "GBS" -> "GESIS" ---
$shp2 = $s.Shapes.Item(2)
$tr2 = $shp2.TextFrame.TextRange
$old2 = $tr2.Characters(1, 3)
$old2.Text = "GESIS"
